# =====================================================================
# Edit script: add a new "Đơn thu nợ" sheet between "Đơn sale chính" and
# "Lương", populate it with a debt-collection order report, and update
# the "Lương" (salary) sheet to add "Chiết khấu thu nợ" (debt-collection
# discount) line items per location plus refreshed totals.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Đơn thu nợ" worksheet right before "Lương" so the
#    tab order becomes: Đơn sale chính, Đơn thu nợ, Lương.
# ---------------------------------------------------------------------
$luongSheet = $wb.Worksheets.Item("Lương")
$newSheet = $wb.Worksheets.Add($luongSheet)
$newSheet.Name = "Đơn thu nợ"

# ---------------------------------------------------------------------
# 2. Populate header row (row 1) of "Đơn thu nợ".
# ---------------------------------------------------------------------
$headers = @(
    "Tiền tố",
    "Mã đơn thu nợ",
    "Lượng thu",
    "Ngày thu",
    "Cơ sở",
    "Đơn nợ",
    "Tên dịch vụ",
    "Khách hàng",
    "Nguồn khách",
    "Sale chính",
    "Đơn giá gốc",
    "Sale phụ",
    "Upsale",
    "Đơn giá",
    "Đã thanh toán",
    "Bác sĩ 1",
    "Bác sĩ 2",
    "Tỉ lệ chiết khấu sale chính",
    "Chiết khấu sale chính",
    "Tỉ lệ chiết khấu sale phụ",
    "Chiết khấu sale phụ",
    "Tỉ lệ chiết khấu bác sĩ 1",
    "Chiết khấu bác sĩ 1",
    "Tỉ lệ chiết khấu bác sĩ 2",
    "Chiết khấu bác sĩ 2"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 3. Row 2: the single debt-collection order.
# ---------------------------------------------------------------------
$newSheet.Cells.Item(2, 1).Value = "TN"
$newSheet.Cells.Item(2, 2).Value = 173
$newSheet.Cells.Item(2, 3).Value = 10000000
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Cells.Item(2, 4).Value = "08-03-2024"
$newSheet.Cells.Item(2, 5).Value = "CẦN THƠ"
$newSheet.Cells.Item(2, 6).Value = "HD-LUXURY-587"
$newSheet.Cells.Item(2, 7).Value = "Nâng mũi"
$newSheet.Cells.Item(2, 8).Value = "Lâm Minh Ngọc"
$newSheet.Cells.Item(2, 9).Value = "Khách cũ"
$newSheet.Cells.Item(2, 10).Value = "Lê Văn Linh"
$newSheet.Cells.Item(2, 11).Value = 20000000
$newSheet.Cells.Item(2, 14).Value = 20000000
$newSheet.Cells.Item(2, 15).Value = 20000000
$newSheet.Cells.Item(2, 16).Value = "Phạm Thanh Hoàng"
$newSheet.Cells.Item(2, 18).Value = 0.1
$newSheet.Cells.Item(2, 19).Value = 1000000
$newSheet.Cells.Item(2, 20).Value = 0
$newSheet.Cells.Item(2, 21).Value = 0
$newSheet.Cells.Item(2, 22).Value = 0
$newSheet.Cells.Item(2, 23).Value = 0
$newSheet.Cells.Item(2, 24).Value = 0
$newSheet.Cells.Item(2, 25).Value = 0

# ---------------------------------------------------------------------
# 4. Row 3: totals row.
# ---------------------------------------------------------------------
$newSheet.Cells.Item(3, 1).Value = "Tổng"
$newSheet.Cells.Item(3, 2).Value = 1
$newSheet.Cells.Item(3, 3).Value = 10000000
$newSheet.Cells.Item(3, 11).Value = 20000000
$newSheet.Cells.Item(3, 13).Value = 0
$newSheet.Cells.Item(3, 14).Value = 20000000
$newSheet.Cells.Item(3, 15).Value = 20000000
$newSheet.Cells.Item(3, 18).Value = 0
$newSheet.Cells.Item(3, 19).Value = 1000000
$newSheet.Cells.Item(3, 20).Value = 0
$newSheet.Cells.Item(3, 21).Value = 0
$newSheet.Cells.Item(3, 22).Value = 0
$newSheet.Cells.Item(3, 23).Value = 0
$newSheet.Cells.Item(3, 24).Value = 0
$newSheet.Cells.Item(3, 25).Value = 0

# ---------------------------------------------------------------------
# 5. Update the "Lương" sheet: insert a "Chiết khấu thu nợ tại <CS>" row
#    right before the "Ứng lương tại <CS>" row of each location.
#    Insert from the bottom up so earlier row numbers stay stable.
# ---------------------------------------------------------------------
$luong = $wb.Worksheets.Item("Lương")

# SÓC TRĂNG block: "Ứng lương tại SÓC TRĂNG" currently at row 31.
$luong.Rows.Item(31).Insert()
$luong.Cells.Item(31, 1).Value = "Chiết khấu thu nợ tại SÓC TRĂNG"
$luong.Cells.Item(31, 2).Value = 0

# LONG XUYÊN block: "Ứng lương tại LONG XUYÊN" currently at row 21.
$luong.Rows.Item(21).Insert()
$luong.Cells.Item(21, 1).Value = "Chiết khấu thu nợ tại LONG XUYÊN"
$luong.Cells.Item(21, 2).Value = 0

# CẦN THƠ block: "Ứng lương tại CẦN THƠ" currently at row 11.
$luong.Rows.Item(11).Insert()
$luong.Cells.Item(11, 1).Value = "Chiết khấu thu nợ tại CẦN THƠ"
$luong.Cells.Item(11, 2).Value = 1000000

# ---------------------------------------------------------------------
# 6. Refresh the grand-total rows (now shifted down by 3 to rows 35-38).
# ---------------------------------------------------------------------
$luong.Cells.Item(35, 1).Value = "Tổng lương tại CẦN THƠ"
$luong.Cells.Item(35, 2).Value = 2141428.571428571
$luong.Cells.Item(36, 1).Value = "Tổng lương tại LONG XUYÊN"
$luong.Cells.Item(36, 2).Value = 714285.7142857143
$luong.Cells.Item(37, 1).Value = "Tổng lương tại SÓC TRĂNG"
$luong.Cells.Item(37, 2).Value = 1071428.571428571
$luong.Cells.Item(38, 1).Value = "Tổng lương tại HỆ THỐNG"
$luong.Cells.Item(38, 2).Value = 3927142.857142857
